# Add data for 2023-12-20
# Updates year-2023 (column J) crime totals across Citywide Totals,
# By Neighborhood, and each affected neighborhood sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 7432
$ws.Range('J3').Value = 7825
$ws.Range('J4').Value = 1704
$ws.Range('J6').Value = 10669
$ws.Range('J7').Value = 28242

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J3').Value = 55
$ws.Range('J7').Value = 424

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 467
$ws.Range('J3').Value = 513
$ws.Range('J6').Value = 659
$ws.Range('J7').Value = 1777

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J6').Value = 152
$ws.Range('J7').Value = 567

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 291
$ws.Range('J3').Value = 427
$ws.Range('J7').Value = 1281

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 261
$ws.Range('J3').Value = 289
$ws.Range('J6').Value = 253
$ws.Range('J7').Value = 866

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 206
$ws.Range('J7').Value = 709

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 167
$ws.Range('J7').Value = 428

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J7').Value = 805
$ws.Range('J8').Value = 1777
$ws.Range('J9').Value = 141
$ws.Range('J11').Value = 505
$ws.Range('J19').Value = 817
$ws.Range('J20').Value = 612
$ws.Range('J23').Value = 260
$ws.Range('J24').Value = 98
$ws.Range('J25').Value = 146
$ws.Range('J29').Value = 1505
$ws.Range('J32').Value = 46
$ws.Range('J33').Value = 1281
$ws.Range('J34').Value = 130
$ws.Range('J37').Value = 866
$ws.Range('J40').Value = 64
$ws.Range('J41').Value = 211
$ws.Range('J44').Value = 222
$ws.Range('J50').Value = 171
$ws.Range('J52').Value = 715
$ws.Range('J53').Value = 424
$ws.Range('J54').Value = 561
$ws.Range('J59').Value = 33
$ws.Range('J63').Value = 86
$ws.Range('J65').Value = 709
$ws.Range('J67').Value = 1031
$ws.Range('J73').Value = 276
$ws.Range('J78').Value = 327
$ws.Range('J79').Value = 772
$ws.Range('J83').Value = 567
$ws.Range('J84').Value = 236
$ws.Range('J85').Value = 1164
$ws.Range('J86').Value = 172
$ws.Range('J88').Value = 299
$ws.Range('J89').Value = 354
$ws.Range('J96').Value = 321
$ws.Range('J99').Value = 428
$ws.Range('J101').Value = 28242

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 390
$ws.Range('J4').Value = 67
$ws.Range('J6').Value = 285
$ws.Range('J7').Value = 1031

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J2').Value = 73
$ws.Range('J7').Value = 236

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 139
$ws.Range('J6').Value = 260
$ws.Range('J7').Value = 561

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 460
$ws.Range('J3').Value = 529
$ws.Range('J7').Value = 1505

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J4').Value = 39
$ws.Range('J6').Value = 317
$ws.Range('J7').Value = 817

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J6').Value = 90
$ws.Range('J7').Value = 222

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J3').Value = 33
$ws.Range('J7').Value = 211

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J2').Value = 86
$ws.Range('J7').Value = 327

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('J3').Value = 27
$ws.Range('J7').Value = 98

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('J2').Value = 72
$ws.Range('J3').Value = 86
$ws.Range('J7').Value = 260

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J3').Value = 82
$ws.Range('J6').Value = 125
$ws.Range('J7').Value = 321

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J2').Value = 220
$ws.Range('J6').Value = 232
$ws.Range('J7').Value = 772

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 170
$ws.Range('J3').Value = 200
$ws.Range('J6').Value = 177
$ws.Range('J7').Value = 612

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J3').Value = 241
$ws.Range('J7').Value = 805

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('J3').Value = 34
$ws.Range('J7').Value = 130

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('J2').Value = 60
$ws.Range('J7').Value = 146

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J4').Value = 27
$ws.Range('J7').Value = 171

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J6').Value = 241
$ws.Range('J7').Value = 505

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('J3').Value = 46
$ws.Range('J7').Value = 141

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J3').Value = 72
$ws.Range('J7').Value = 276

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range('J6').Value = 8
$ws.Range('J7').Value = 33

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J2').Value = 61
$ws.Range('J7').Value = 299

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range('J2').Value = 15
$ws.Range('J7').Value = 46

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J6').Value = 110
$ws.Range('J7').Value = 354

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('J4').Value = 93
$ws.Range('J7').Value = 172

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 310
$ws.Range('J3').Value = 420
$ws.Range('J6').Value = 333
$ws.Range('J7').Value = 1164

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('J3').Value = 24
$ws.Range('J7').Value = 64

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J4').Value = 27
$ws.Range('J7').Value = 715
